$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet ---
$wsSchedule.Range("B2").Value = 46043.1875
$wsSchedule.Range("C2").Value = 4.5
$wsSchedule.Range("D2").Value = 17.01
$wsSchedule.Range("E2").Value = 573.8929364999999
$wsSchedule.Range("F2").Value = 33.73856181657848
$wsSchedule.Range("A4").Value = 46043.89583333334
$wsSchedule.Range("B4").Value = 46044.10416666666
$wsSchedule.Range("C4").Value = 5
$wsSchedule.Range("D4").Value = 18.9
$wsSchedule.Range("E4").Value = 588.914547
$wsSchedule.Range("F4").Value = 31.15949984126984
$wsSchedule.Range("A5").Value = 46044.27083333334
$wsSchedule.Range("C5").Value = 9.5
$wsSchedule.Range("D5").Value = 35.91
$wsSchedule.Range("E5").Value = 13.66849575
$wsSchedule.Range("F5").Value = 0.3806320175438597
# --- Detailed sheet ---
$wsDetailed.Range("E10").Value = "ON"
$wsDetailed.Range("B35").Value = -5.51
$wsDetailed.Range("B36").Value = 36.06
$wsDetailed.Range("B37").Value = 36.08141
$wsDetailed.Range("B38").Value = 46.54611
$wsDetailed.Range("C38").Value = "historical"
$wsDetailed.Range("B39").Value = 64.35683
$wsDetailed.Range("B40").Value = 73.19
$wsDetailed.Range("B41").Value = 77.94
$wsDetailed.Range("B42").Value = 100.01
$wsDetailed.Range("B43").Value = 74.43841999999999
$wsDetailed.Range("B44").Value = 73.19
$wsDetailed.Range("E44").Value = "OFF"
$wsDetailed.Range("B45").Value = 71.7281
$wsDetailed.Range("B46").Value = 59.62291
$wsDetailed.Range("B47").Value = 57.98348
$wsDetailed.Range("B48").Value = 60.73801
$wsDetailed.Range("B49").Value = 61.32156
$wsDetailed.Range("B51").Value = 63.88086
$wsDetailed.Range("E55").Value = "OFF"
$wsDetailed.Range("B56").Value = 63.73519
$wsDetailed.Range("B57").Value = 64.10364
$wsDetailed.Range("B58").Value = 65.32088
$wsDetailed.Range("B59").Value = 65.85026999999999
$wsDetailed.Range("B60").Value = 66.04559
$wsDetailed.Range("B61").Value = 76.12006
$wsDetailed.Range("B62").Value = 64.89
$wsDetailed.Range("E63").Value = "ON"
$wsDetailed.Range("B64").Value = 30.9379
$wsDetailed.Range("B66").Value = -5.50985
$wsDetailed.Range("B67").Value = -6.15086
$wsDetailed.Range("B68").Value = -10
$wsDetailed.Range("B69").Value = -13.50737
$wsDetailed.Range("B70").Value = -9.621499999999999
$wsDetailed.Range("B71").Value = -12.01
$wsDetailed.Range("B72").Value = -12.01
$wsDetailed.Range("B73").Value = -5.74313
$wsDetailed.Range("B74").Value = -10
$wsDetailed.Range("B75").Value = -8.0564
$wsDetailed.Range("B76").Value = -7.79393
$wsDetailed.Range("B77").Value = -5.88864
$wsDetailed.Range("B78").Value = -5.27725
$wsDetailed.Range("B79").Value = 0.51
$wsDetailed.Range("B81").Value = 36.06
$wsDetailed.Range("B82").Value = 0.51
$wsDetailed.Range("B83").Value = -4.13512
$wsDetailed.Range("B84").Value = -5.14805
$wsDetailed.Range("B85").Value = -6.90848
$wsDetailed.Range("B86").Value = -0.45834
$wsDetailed.Range("B87").Value = 0.00036
$wsDetailed.Range("B88").Value = 10.48193
$wsDetailed.Range("B89").Value = 55.33036
$wsDetailed.Range("B90").Value = 53.90468
$wsDetailed.Range("B91").Value = 54.47327
$wsDetailed.Range("B92").Value = 57.01318
$wsDetailed.Range("B93").Value = 50.38252
$wsDetailed.Range("B94").Value = 30.67112
$wsDetailed.Range("B95").Value = 56.98
$wsDetailed.Range("B96").Value = 56.23018
$wsDetailed.Range("B97").Value = 48.31676
